# Automatische test-sync: 2025-06-22 17:34:50
# Adds the new "Afmelding nieuwsbrief" mail-log entry (Logs row 11) and its
# corresponding Dashboard roll-up row (Dashboard row 9), then widens the
# conditional-formatting ranges and chart series to cover the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append the new row of data (row 11)
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Afmelding nieuwsbrief"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D11").Value = "Afmelding / Nieuwsbrief"
$logs.Range("E11").Value = "Geachte heer/mevrouw,`nHartelijk dank voor uw e-mail. Ik heb u bij deze afgemeld voor de nieuwsbrief. Mocht u in de toekomst weer interesse hebben, dan kunt u zich altijd opnieuw aanmelden.`nMet vriendelijke groet,`n[Naam]"
$logs.Range("F11").Value = "2025-06-22 17:34:14"
$logs.Range("G11").Value = "Ja"

# The multi-line answer in E11 would otherwise pin an explicit (and noisy)
# custom row height; auto-fit it back down so row 11 matches the plain,
# unsized rows above it.
$logs.Rows.Item(11).AutoFit()

# Widen the conditional formatting ranges so the new row is covered too.
$logs.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: append the roll-up row (row 9) for the new category
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A9").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B9").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the category/value series refs
#    from row 8 to row 9 so the new category shows up in the chart too.
# ---------------------------------------------------------------------------
$cho = $dash.ChartObjects(1)
$chart = $cho.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$9,'Dashboard'!`$B`$2:`$B`$9,1)"
